$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- data from former row 16
$ws.Range("D2").Value = 45021
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 4500
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 4700
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 783

# Row 3 <- data from former row 15
$ws.Range("D3").Value = 44987
$ws.Range("J3").Value = 130
$ws.Range("K3").Value = 4500
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = 4692
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 782

# Row 4 <- data from former row 10
$ws.Range("D4").Value = 44650
$ws.Range("J4").Value = 130
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3500
$ws.Range("M4").Value = 3308
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 551

# Row 5 <- data from former row 13
$ws.Range("D5").Value = 44631
$ws.Range("J5").Value = 110
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3273
$ws.Range("O5").Value = "Provincia de Chacabuco"
$ws.Range("P5").Value = 546

# Row 6 <- data from former row 9
$ws.Range("D6").Value = 44672
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3500
$ws.Range("M6").Value = 3286
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 548

# Row 7 <- data from former row 4
$ws.Range("D7").Value = 44630
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 2722
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 454

# Row 8 <- data from former row 12
$ws.Range("D8").Value = 44685
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3500
$ws.Range("M8").Value = 3267
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 544

# Row 9 <- data from former row 6
$ws.Range("D9").Value = 44658
$ws.Range("J9").Value = 180
$ws.Range("K9").Value = 2500
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 2778
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 463

# Row 10 <- data from former row 8
$ws.Range("D10").Value = 44643
$ws.Range("J10").Value = 90
$ws.Range("K10").Value = 2800
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2911
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 485

# Row 12 <- data from former row 2
$ws.Range("D12").Value = 44644
$ws.Range("J12").Value = 140
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 2786
$ws.Range("O12").Value = "Provincia de Chacabuco"
$ws.Range("P12").Value = 464

# Row 13 <- data from former row 7
$ws.Range("D13").Value = 44876
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 6500
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6812
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 1135

# Row 14 <- data from former row 3
$ws.Range("D14").Value = 44957
$ws.Range("J14").Value = 70
$ws.Range("K14").Value = 1500
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 1857
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 310

# Row 15 <- data from former row 14
$ws.Range("D15").Value = 44637
$ws.Range("J15").Value = 170
$ws.Range("K15").Value = 2800
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 2906
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 484

# Row 16 <- data from former row 5
$ws.Range("D16").Value = 44659
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 2722
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 454
